$d = $word.ActiveDocument

$d.Range(3452, 3645).Text = '5817692 - Katia Cristiane Gandolpho Candioto'
$d.Range(2638, 3438).Text = '519033 - Carlos Yujiro Shigue'
$d.Range(2465, 2616).Text = 'Os indicadores serão obtidos por questionário de avaliação pelos usuários quanto aos seguintes quesitos: conhecimento adquirido e satisfação do usuário nas apresentações e formas de divulgação.' + [char]11 + ''
$d.Range(2402, 2455).Text = '- Identificação das necessidades do grupo social: pesquisas, entrevistas e observações para entender as necessidades, desafios e preferências dos estudantes.' + [char]11 + '- Definição de objetivos e requisitos do projeto para que as soluções desenvolvidas devem atender: identificar funcionalidades, restrições de orçamento e cronograma, e quaisquer outras considerações importantes.' + [char]11 + '- Pesquisa e desenvolvimento projetos relacionados à engenharia: criação de protótipos, desenvolvimento de software, fabricação de dispositivos e apresentação de aplicações para garantir que haja disseminação do conhecimento sobre a profissão engenharia.' + [char]11 + '- Avaliação: feedback recebido quanto ao conhecimento sobre o tema.' + [char]11 + '- Implementação e distribuição: Visita e apresentações em escolas de ensino infantil, fundamental ou médio.' + [char]11 + ''
$d.Range(1514, 1906).Text = 'Para os estudantes: despertar interesse na engenharia.' + [char]11 + 'Para a formação dos discentes: Desenvolver conceitos de engenharia com aplicações profissionais'
$d.Range(1437, 1504).Text = 'Introduce students to the principles and methodology of scientific research.'
$d.Range(1356, 1436).Text = 'Estudantes de ensino infantil, fundamental ou médio.'
$d.Range(1288, 1337).Text = 'Nota de avaliação do projeto e demais documentos.' + [char]11 + 'Devido às características práticas da disciplina, não será oferecida recuperação' + [char]11 + 'ASTI VERA, A. Metodologia da pesquisa científica. Porto Alegre: Ed. Globo, 1973. BARRAS, R. Os cientistas precisam escrever: guia de redação para cientistas, engenheiros e estudantes. São Paulo: TAQ/EDUSP, 1979. CERVO, A. L.; BERVIAN, P. A. Metodologia científica. São Paulo: Mc-Graw-Hill do Brasil, 1973. ANDRADE, M. M. Introdução à Metodologia do Trabalho Científico São Paulo: Atlas, 2005.'
$d.Range(479, 546).Text = 'Introduzir aos estudantes os princípios e a metodologia da pesquisa científica.' + [char]11 + ''
$d.Range(404, 479).Text = ''
$d.Range(324, 404).Text = 'Iniciação a um projeto de pesquisa sob orientação de um professor.' + [char]11 + ''
$d.Range(218, 294).Text = 'Initiation to a research project under the guidance of a professor.'
